# Fix Training Data Issue (#48)
# ------------------------------------------------------------------
# The "Date" column was populated from the source filename
# (6-23-2011-12), which is one day off from the actual NBA game date
# because of the way NBA.com displayed stats. Correct it to the real
# ISO date, 2012-06-23, for every data row on the sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "6-23-2011-12"
$newValue = "2012-06-23"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol = $firstCol + $used.Columns.Count - 1

# Locate the "Date" header column on the first row so the fix keeps
# working even if the sheet layout shifts.
$dateCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Value2
    if ($header -eq "Date") {
        $dateCol = $c
    }
}

if ($dateCol -eq 0) {
    $dateCol = 58 # fall back to column BF
}

for ($row = $firstRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq $oldValue) {
        # Assign with a leading apostrophe so the engine keeps storing a
        # plain text value instead of re-parsing "2012-06-23" as a date
        # serial number, then drop back to the workbook's default style
        # so no extra number-format gets attached to the cell.
        $cell.Value = "'" + $newValue
        $cell.Style = "Normal"
    }
}
